# Update marksheet correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: "Marking" - correct answer marking value B11: 3 -> 5
$ws.Range("B11").Value = 5

# Row 12: "Total" - total marks B12: 54 -> 90
$ws.Range("B12").Value = 90

# Row 12: "Total" correct/total fraction text E12: "51/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
